$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet: last-updated timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:06 PM"

# --- Update Stock List sheet: new data pushes two new rows to the top ---
$ws = $wb.Worksheets.Item("Stock List")

# Insert two blank rows at the top of the data (row 2), shifting existing
# data down by two rows.
$ws.Range("A2:A3").EntireRow.Insert()

# The two rows that fell off the bottom of the list (now at 77:78) are removed.
$ws.Range("A77:A78").EntireRow.Delete()

# Clear any inherited formatting on the two new rows so they match the
# plain (unstyled) data rows.
$ws.Range("A2:H3").ClearFormats()

# Fill in the new row 2 (MIDWESTLTD)
$ws.Range("A2").Value = "📋"
$ws.Range("B2").Value = "MIDWESTLTD"
$ws.Range("C2").Value = "MIDWESTLTD"
$ws.Range("D2").Value = 1117.2
$ws.Range("E2").Value = -1.4032
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = 4039.8864

# Fill in the new row 3 (CAPTRU-RE1)
$ws.Range("A3").Value = "📋"
$ws.Range("B3").Value = "CAPTRU-RE1"
$ws.Range("C3").Value = "CAPTRU-RE1"
$ws.Range("D3").Value = 5.67
$ws.Range("E3").Value = -11.9565
$ws.Range("F3").Value = "N/A"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = 0
